{"js": "// Auto-generated edit: replace 100 arithmetic expressions in the single\n// table (20 rows x 5 cols) with updated expressions, in row-major order,\n// matching the order the cells appear in the document.\nconst replacements = [[\"39-6=\", \"13-6=\"], [\"85-63=\", \"79-58=\"], [\"3+92=\", \"10+87=\"], [\"89-78=\", \"91-44=\"], [\"0+92=\", \"34+43=\"], [\"91-62=\", \"92-20=\"], [\"12+16=\", \"31+34=\"], [\"38-26=\", \"81-39=\"], [\"54-0=\", \"90-54=\"], [\"37+20=\", \"58-26=\"], [\"19+65=\", \"15+73=\"], [\"47+39=\", \"20+33=\"], [\"1+51=\", \"51+1=\"], [\"94-71=\", \"73-38=\"], [\"29+53=\", \"14+2=\"], [\"23-0=\", \"88-9=\"], [\"74-36=\", \"65+18=\"], [\"35-2=\", \"69-15=\"], [\"73-48=\", \"46-3=\"], [\"98-37=\", \"38+31=\"], [\"67-61=\", \"6+60=\"], [\"67-62=\", \"77-61=\"], [\"89+8=\", \"81-33=\"], [\"1+24=\", \"35+49=\"], [\"67-55=\", \"10+33=\"], [\"35-12=\", \"91-0=\"], [\"47+45=\", \"67-15=\"], [\"23+4=\", \"38-28=\"], [\"55-51=\", \"18-18=\"], [\"66-35=\", \"86-33=\"], [\"31+35=\", \"29+70=\"], [\"13+13=\", \"68+7=\"], [\"57+18=\", \"32-12=\"], [\"95-6=\", \"49+46=\"], [\"18+32=\", \"71-9=\"], [\"57-14=\", \"19+4=\"], [\"77+4=\", \"90-13=\"], [\"20+50=\", \"56+43=\"], [\"42+25=\", \"95-90=\"], [\"17+59=\", \"78+13=\"], [\"0+64=\", \"18+66=\"], [\"46+10=\", \"16+56=\"], [\"95-18=\", \"26-21=\"], [\"87-29=\", \"23+34=\"], [\"59-32=\", \"69+18=\"], [\"76+3=\", \"73-11=\"], [\"35-34=\", \"33+61=\"], [\"49-47=\", \"13+22=\"], [\"88-20=\", \"88-3=\"], [\"36-31=\", \"26+53=\"], [\"29-23=\", \"14+75=\"], [\"35-27=\", \"69-6=\"], [\"6+42=\", \"1+74=\"], [\"60-31=\", \"15+73=\"], [\"16+37=\", \"89-17=\"], [\"75-20=\", \"95-3=\"], [\"3+84=\", \"45+2=\"], [\"64-33=\", \"99-5=\"], [\"7+54=\", \"39+52=\"], [\"37+42=\", \"14+15=\"], [\"73+13=\", \"30+40=\"], [\"94-63=\", \"63+33=\"], [\"33+32=\", \"56+14=\"], [\"27-4=\", \"66-46=\"], [\"34-9=\", \"59-0=\"], [\"54-7=\", \"83-55=\"], [\"49+11=\", \"69-15=\"], [\"89-36=\", \"97-45=\"], [\"21-19=\", \"71+0=\"], [\"20+41=\", \"49-7=\"], [\"9+51=\", \"24+70=\"], [\"79-39=\", \"87+3=\"], [\"2+65=\", \"69+5=\"], [\"11+45=\", \"20+11=\"], [\"20-0=\", \"0+4=\"], [\"30-14=\", \"33-14=\"], [\"97-72=\", \"40-12=\"], [\"85+2=\", \"40-38=\"], [\"77-74=\", \"40+2=\"], [\"92-78=\", \"41+7=\"], [\"29+50=\", \"8+7=\"], [\"24+51=\", \"87-28=\"], [\"72-3=\", \"39+2=\"], [\"60-56=\", \"17+73=\"], [\"84-23=\", \"17-8=\"], [\"62-35=\", \"22+32=\"], [\"76-27=\", \"90-48=\"], [\"28-4=\", \"11+28=\"], [\"61-27=\", \"21+55=\"], [\"13+19=\", \"20+31=\"], [\"96+1=\", \"16+19=\"], [\"42+52=\", \"19+66=\"], [\"32+34=\", \"13+14=\"], [\"44-2=\", \"40+2=\"], [\"83-48=\", \"8+57=\"], [\"22+63=\", \"70-65=\"], [\"64+34=\", \"92-29=\"], [\"52-10=\", \"3+78=\"], [\"15+65=\", \"91-22=\"], [\"52-11=\", \"76-50=\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document\");\n}\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nconst firstRow = table.rows.getFirst();\nfirstRow.load(\"cellCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = firstRow.cellCount;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    if (idx >= replacements.length) break;\n    const [oldText, newText] = replacements[idx];\n    const cell = table.getCell(r, c);\n    const par = cell.body.paragraphs.getFirst();\n    const range = par.getRange();\n    range.load(\"text\");\n    await context.sync();\n    if (range.text !== oldText) {\n      throw new Error(\n        \"Mismatch at row \" + r + \" col \" + c + \": expected '\" + oldText +\n        \"' but found '\" + range.text + \"'\"\n      );\n    }\n    range.insertText(newText, Word.InsertLocation.replace);\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Auto-generated edit: replace 100 arithmetic expressions in the single\n# table (20 rows x 5 cols) with updated expressions, in row-major order,\n# matching the order the cells appear in the document.\n$replacements = @(\n    @(\"39-6=\",\"13-6=\"),\n    @(\"85-63=\",\"79-58=\"),\n    @(\"3+92=\",\"10+87=\"),\n    @(\"89-78=\",\"91-44=\"),\n    @(\"0+92=\",\"34+43=\"),\n    @(\"91-62=\",\"92-20=\"),\n    @(\"12+16=\",\"31+34=\"),\n    @(\"38-26=\",\"81-39=\"),\n    @(\"54-0=\",\"90-54=\"),\n    @(\"37+20=\",\"58-26=\"),\n    @(\"19+65=\",\"15+73=\"),\n    @(\"47+39=\",\"20+33=\"),\n    @(\"1+51=\",\"51+1=\"),\n    @(\"94-71=\",\"73-38=\"),\n    @(\"29+53=\",\"14+2=\"),\n    @(\"23-0=\",\"88-9=\"),\n    @(\"74-36=\",\"65+18=\"),\n    @(\"35-2=\",\"69-15=\"),\n    @(\"73-48=\",\"46-3=\"),\n    @(\"98-37=\",\"38+31=\"),\n    @(\"67-61=\",\"6+60=\"),\n    @(\"67-62=\",\"77-61=\"),\n    @(\"89+8=\",\"81-33=\"),\n    @(\"1+24=\",\"35+49=\"),\n    @(\"67-55=\",\"10+33=\"),\n    @(\"35-12=\",\"91-0=\"),\n    @(\"47+45=\",\"67-15=\"),\n    @(\"23+4=\",\"38-28=\"),\n    @(\"55-51=\",\"18-18=\"),\n    @(\"66-35=\",\"86-33=\"),\n    @(\"31+35=\",\"29+70=\"),\n    @(\"13+13=\",\"68+7=\"),\n    @(\"57+18=\",\"32-12=\"),\n    @(\"95-6=\",\"49+46=\"),\n    @(\"18+32=\",\"71-9=\"),\n    @(\"57-14=\",\"19+4=\"),\n    @(\"77+4=\",\"90-13=\"),\n    @(\"20+50=\",\"56+43=\"),\n    @(\"42+25=\",\"95-90=\"),\n    @(\"17+59=\",\"78+13=\"),\n    @(\"0+64=\",\"18+66=\"),\n    @(\"46+10=\",\"16+56=\"),\n    @(\"95-18=\",\"26-21=\"),\n    @(\"87-29=\",\"23+34=\"),\n    @(\"59-32=\",\"69+18=\"),\n    @(\"76+3=\",\"73-11=\"),\n    @(\"35-34=\",\"33+61=\"),\n    @(\"49-47=\",\"13+22=\"),\n    @(\"88-20=\",\"88-3=\"),\n    @(\"36-31=\",\"26+53=\"),\n    @(\"29-23=\",\"14+75=\"),\n    @(\"35-27=\",\"69-6=\"),\n    @(\"6+42=\",\"1+74=\"),\n    @(\"60-31=\",\"15+73=\"),\n    @(\"16+37=\",\"89-17=\"),\n    @(\"75-20=\",\"95-3=\"),\n    @(\"3+84=\",\"45+2=\"),\n    @(\"64-33=\",\"99-5=\"),\n    @(\"7+54=\",\"39+52=\"),\n    @(\"37+42=\",\"14+15=\"),\n    @(\"73+13=\",\"30+40=\"),\n    @(\"94-63=\",\"63+33=\"),\n    @(\"33+32=\",\"56+14=\"),\n    @(\"27-4=\",\"66-46=\"),\n    @(\"34-9=\",\"59-0=\"),\n    @(\"54-7=\",\"83-55=\"),\n    @(\"49+11=\",\"69-15=\"),\n    @(\"89-36=\",\"97-45=\"),\n    @(\"21-19=\",\"71+0=\"),\n    @(\"20+41=\",\"49-7=\"),\n    @(\"9+51=\",\"24+70=\"),\n    @(\"79-39=\",\"87+3=\"),\n    @(\"2+65=\",\"69+5=\"),\n    @(\"11+45=\",\"20+11=\"),\n    @(\"20-0=\",\"0+4=\"),\n    @(\"30-14=\",\"33-14=\"),\n    @(\"97-72=\",\"40-12=\"),\n    @(\"85+2=\",\"40-38=\"),\n    @(\"77-74=\",\"40+2=\"),\n    @(\"92-78=\",\"41+7=\"),\n    @(\"29+50=\",\"8+7=\"),\n    @(\"24+51=\",\"87-28=\"),\n    @(\"72-3=\",\"39+2=\"),\n    @(\"60-56=\",\"17+73=\"),\n    @(\"84-23=\",\"17-8=\"),\n    @(\"62-35=\",\"22+32=\"),\n    @(\"76-27=\",\"90-48=\"),\n    @(\"28-4=\",\"11+28=\"),\n    @(\"61-27=\",\"21+55=\"),\n    @(\"13+19=\",\"20+31=\"),\n    @(\"96+1=\",\"16+19=\"),\n    @(\"42+52=\",\"19+66=\"),\n    @(\"32+34=\",\"13+14=\"),\n    @(\"44-2=\",\"40+2=\"),\n    @(\"83-48=\",\"8+57=\"),\n    @(\"22+63=\",\"70-65=\"),\n    @(\"64+34=\",\"92-29=\"),\n    @(\"52-10=\",\"3+78=\"),\n    @(\"15+65=\",\"91-22=\"),\n    @(\"52-11=\",\"76-50=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $replacements.Count) { break }\n        $pair = $replacements[$idx]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n        $cell = $t.Cell($r, $c)\n        $range = $cell.Range\n        $current = $range.Text.TrimEnd([char]13, [char]7)\n        if ($current -ne $oldText) {\n            throw \"Mismatch at row $r col $c`: expected '$oldText' but found '$current'\"\n        }\n        $range.Text = $newText\n        $idx++\n    }\n}\n"}
